$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").EntireColumn.Insert()
$d7 = $ws.Range("D7").Value2
$e7 = $ws.Range("E7").Value2
"D7: $d7"
"E7: $e7"
